$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 (shifts existing rows 70-123 down to 71-124)
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new record's data
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value = "Bíobío"
$ws.Cells.Item(70, 4).Value = 44651
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100109
$ws.Cells.Item(70, 8).Value = "Uva"
$ws.Cells.Item(70, 9).Value = 100109001
$ws.Cells.Item(70, 10).Value = "Uva"
$ws.Cells.Item(70, 11).Value = "Thompson seedless"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 150
$ws.Cells.Item(70, 14).Value = 10000
$ws.Cells.Item(70, 15).Value = 11000
$ws.Cells.Item(70, 16).Value = 10533
$ws.Cells.Item(70, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(70, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 19).Value = 585
$ws.Cells.Item(70, 20).Value = 18
